# Timing for Matrix Factorization
#
# - user_based / item_based: add an "Average Runtime" column (C) and a new
#   11-13-2018 data row, and give the two existing date rows real date
#   values (formatted m-d-yyyy) instead of the old blank/text placeholders.
# - Add a new Matrix_Factorization sheet with the same Date / Best Runtime /
#   Average Runtime layout as the other two timing sheets.

$wb = $excel.ActiveWorkbook

function Fill-TimingSheet($ws, $bestDay1, $bestDay2, $avgDay3) {
    # --- Headers (row 1) ---
    $ws.Range("A1").Value = "Date"
    $ws.Range("B1").Value = "Best Runtime"

    # New "Average Runtime" header - reuse the same (gray-filled) style as
    # the other header cells.
    $ws.Range("A1").Copy() | Out-Null
    $ws.Range("C1").PasteSpecial(-4122) | Out-Null
    $ws.Range("C1").Value = "Average Runtime"

    # --- Row 2: 11-11-2018 ---
    $ws.Range("A2").NumberFormat = "m-d-yyyy"
    $ws.Range("A2").Value = 43415
    $ws.Range("B2").Value = $bestDay1

    # --- Row 3: 11-12-2018 ---
    $ws.Range("A3").NumberFormat = "m-d-yyyy"
    $ws.Range("A3").Value = 43416
    $ws.Range("B3").Value = $bestDay2

    # --- Row 4: 11-13-2018, stored as literal text like the original row ---
    $ws.Range("A4").Formula = '="11-13-2018"'
    $ws.Range("A4").Copy() | Out-Null
    $ws.Range("A4").PasteSpecial(-4163) | Out-Null
    $ws.Range("B4").Value = $avgDay3
}

# --- user_based ---
$wsUser = $wb.Worksheets.Item("user_based")
Fill-TimingSheet $wsUser 1.59 1.58 1.54

# --- item_based ---
$wsItem = $wb.Worksheets.Item("item_based")
Fill-TimingSheet $wsItem 4.6 4.57 4.53

# --- new Matrix_Factorization sheet, appended after item_based ---
# (only 2 data rows - no 11-13-2018 row yet for this algorithm)
$wsMF = $wb.Worksheets.Add($null, $wsItem)
$wsMF.Name = "Matrix_Factorization"

# Match the page setup (and thus the default row/col sizing) of the other
# timing sheets: 0.75in left/right, 1in top/bottom, 0.5in header/footer.
$wsMF.PageSetup.LeftMargin = 54
$wsMF.PageSetup.RightMargin = 54
$wsMF.PageSetup.TopMargin = 72
$wsMF.PageSetup.BottomMargin = 72
$wsMF.PageSetup.HeaderMargin = 36
$wsMF.PageSetup.FooterMargin = 36

$wsUser.Range("A1:B1").Copy() | Out-Null
$wsMF.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$wsUser.Range("C1").Copy() | Out-Null
$wsMF.Range("C1").PasteSpecial(-4122) | Out-Null

$wsMF.Range("A1").Value = "Date"
$wsMF.Range("B1").Value = "Best Runtime"
$wsMF.Range("C1").Value = "Average Runtime"

$wsMF.Range("A2").NumberFormat = "m-d-yyyy"
$wsMF.Range("A2").Value = 43415
$wsMF.Range("B2").Value = 1.59

$wsMF.Range("A3").NumberFormat = "m-d-yyyy"
$wsMF.Range("A3").Value = 43416
$wsMF.Range("B3").Value = 1.58
